# FRED WALCL data refresh:
#  - Append the two newest weekly observations to the "Data" sheet
#  - Update the "SeriesInfo" sheet metadata to reflect the new pull

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Data")
$infoSheet = $wb.Worksheets.Item("SeriesInfo")

# --- Data sheet: append rows 110 and 111, matching the date-column
#     formatting already used by the existing observation rows. ---
$dataSheet.Range("A109:B109").Copy() | Out-Null
$dataSheet.Range("A110:B110").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$dataSheet.Range("A109:B109").Copy() | Out-Null
$dataSheet.Range("A111:B111").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = $false

$dataSheet.Cells.Item(110, 1).Value = 45231
$dataSheet.Cells.Item(110, 2).Value = 7866.664

$dataSheet.Cells.Item(111, 1).Value = 45238
$dataSheet.Cells.Item(111, 2).Value = 7860.691

# --- SeriesInfo sheet: refresh the metadata that changes on every pull.
#     Force plain text (not an auto-detected date/number) the same way the
#     existing cells are stored, then drop the scratch formatting again so
#     no stray cell format is left behind. ---
function Set-PlainTextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-PlainTextValue $infoSheet.Range("B3")  "2023-11-15"              # realtime_start
Set-PlainTextValue $infoSheet.Range("B4")  "2023-11-15"              # realtime_end
Set-PlainTextValue $infoSheet.Range("B7")  "2023-11-08"              # observation_end
Set-PlainTextValue $infoSheet.Range("B14") "2023-11-09 15:37:01-06"  # last_updated
